# Update the 25 two-digit multiplication problems that appear in the
# first (and only) table of the document. The problems sit in table
# rows 1, 5, 10, 15 and 20 (the other rows are blank answer rows).
#
# Each cell is targeted explicitly by (row, column) via Tables(1).Cell(),
# and the replacement is scoped to that cell's own Range (rebuilt with
# Document.Range(start, end) to dodge a COM-shim quirk where calling
# Find/Replace straight off a TableCell.Range can bleed into another
# cell that happens to hold identical text, e.g. the two "51x51=" cells).

$d = $word.ActiveDocument
$t = $d.Tables(1)

$cell = $t.Cell(1, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("67×38=", $true, $false, $false, $false, $false, $true, 0, $false, "77×69=", 2)
if (-not $found) { Write-Output "FAILED at row 1 col 1: 67×38= -> 77×69=" }
$cell = $t.Cell(1, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("51×57=", $true, $false, $false, $false, $false, $true, 0, $false, "32×11=", 2)
if (-not $found) { Write-Output "FAILED at row 1 col 2: 51×57= -> 32×11=" }
$cell = $t.Cell(1, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("15×41=", $true, $false, $false, $false, $false, $true, 0, $false, "51×97=", 2)
if (-not $found) { Write-Output "FAILED at row 1 col 3: 15×41= -> 51×97=" }
$cell = $t.Cell(1, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("41×25=", $true, $false, $false, $false, $false, $true, 0, $false, "68×86=", 2)
if (-not $found) { Write-Output "FAILED at row 1 col 4: 41×25= -> 68×86=" }
$cell = $t.Cell(1, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("36×45=", $true, $false, $false, $false, $false, $true, 0, $false, "17×38=", 2)
if (-not $found) { Write-Output "FAILED at row 1 col 5: 36×45= -> 17×38=" }
$cell = $t.Cell(5, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("19×84=", $true, $false, $false, $false, $false, $true, 0, $false, "23×57=", 2)
if (-not $found) { Write-Output "FAILED at row 5 col 1: 19×84= -> 23×57=" }
$cell = $t.Cell(5, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("99×23=", $true, $false, $false, $false, $false, $true, 0, $false, "47×74=", 2)
if (-not $found) { Write-Output "FAILED at row 5 col 2: 99×23= -> 47×74=" }
$cell = $t.Cell(5, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("67×31=", $true, $false, $false, $false, $false, $true, 0, $false, "51×47=", 2)
if (-not $found) { Write-Output "FAILED at row 5 col 3: 67×31= -> 51×47=" }
$cell = $t.Cell(5, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("47×58=", $true, $false, $false, $false, $false, $true, 0, $false, "15×39=", 2)
if (-not $found) { Write-Output "FAILED at row 5 col 4: 47×58= -> 15×39=" }
$cell = $t.Cell(5, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("95×95=", $true, $false, $false, $false, $false, $true, 0, $false, "34×33=", 2)
if (-not $found) { Write-Output "FAILED at row 5 col 5: 95×95= -> 34×33=" }
$cell = $t.Cell(10, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("56×15=", $true, $false, $false, $false, $false, $true, 0, $false, "43×32=", 2)
if (-not $found) { Write-Output "FAILED at row 10 col 1: 56×15= -> 43×32=" }
$cell = $t.Cell(10, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("99×30=", $true, $false, $false, $false, $false, $true, 0, $false, "48×14=", 2)
if (-not $found) { Write-Output "FAILED at row 10 col 2: 99×30= -> 48×14=" }
$cell = $t.Cell(10, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("74×91=", $true, $false, $false, $false, $false, $true, 0, $false, "49×47=", 2)
if (-not $found) { Write-Output "FAILED at row 10 col 3: 74×91= -> 49×47=" }
$cell = $t.Cell(10, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("52×55=", $true, $false, $false, $false, $false, $true, 0, $false, "46×64=", 2)
if (-not $found) { Write-Output "FAILED at row 10 col 4: 52×55= -> 46×64=" }
$cell = $t.Cell(10, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("51×51=", $true, $false, $false, $false, $false, $true, 0, $false, "48×57=", 2)
if (-not $found) { Write-Output "FAILED at row 10 col 5: 51×51= -> 48×57=" }
$cell = $t.Cell(15, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("58×22=", $true, $false, $false, $false, $false, $true, 0, $false, "71×14=", 2)
if (-not $found) { Write-Output "FAILED at row 15 col 1: 58×22= -> 71×14=" }
$cell = $t.Cell(15, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("47×85=", $true, $false, $false, $false, $false, $true, 0, $false, "67×12=", 2)
if (-not $found) { Write-Output "FAILED at row 15 col 2: 47×85= -> 67×12=" }
$cell = $t.Cell(15, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("81×51=", $true, $false, $false, $false, $false, $true, 0, $false, "82×15=", 2)
if (-not $found) { Write-Output "FAILED at row 15 col 3: 81×51= -> 82×15=" }
$cell = $t.Cell(15, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("51×55=", $true, $false, $false, $false, $false, $true, 0, $false, "92×57=", 2)
if (-not $found) { Write-Output "FAILED at row 15 col 4: 51×55= -> 92×57=" }
$cell = $t.Cell(15, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("51×51=", $true, $false, $false, $false, $false, $true, 0, $false, "89×55=", 2)
if (-not $found) { Write-Output "FAILED at row 15 col 5: 51×51= -> 89×55=" }
$cell = $t.Cell(20, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("79×60=", $true, $false, $false, $false, $false, $true, 0, $false, "13×54=", 2)
if (-not $found) { Write-Output "FAILED at row 20 col 1: 79×60= -> 13×54=" }
$cell = $t.Cell(20, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("87×65=", $true, $false, $false, $false, $false, $true, 0, $false, "42×11=", 2)
if (-not $found) { Write-Output "FAILED at row 20 col 2: 87×65= -> 42×11=" }
$cell = $t.Cell(20, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("15×88=", $true, $false, $false, $false, $false, $true, 0, $false, "32×89=", 2)
if (-not $found) { Write-Output "FAILED at row 20 col 3: 15×88= -> 32×89=" }
$cell = $t.Cell(20, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("72×85=", $true, $false, $false, $false, $false, $true, 0, $false, "17×35=", 2)
if (-not $found) { Write-Output "FAILED at row 20 col 4: 72×85= -> 17×35=" }
$cell = $t.Cell(20, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$found = $rng.Find.Execute("61×24=", $true, $false, $false, $false, $false, $true, 0, $false, "61×79=", 2)
if (-not $found) { Write-Output "FAILED at row 20 col 5: 61×24= -> 61×79=" }

Write-Output "Done"
